$d = $word.ActiveDocument

# Locate the "LOQ4057: Operações Unitárias III (Requisito fraco)" paragraph,
# which must be kept. The three paragraphs right after it - a blank
# paragraph, the "Ver no Jupiter Salvar em pdf Salvar em docx" line, and the
# copyright/footer line - are the ones being removed by this edit.
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOQ4057*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate the LOQ4057 paragraph"
}

# Delete from the highest paragraph index down to the lowest so earlier
# deletions don't shift the indices of paragraphs still to be removed.
$d.Paragraphs.Item($anchorIndex + 3).Range.Delete()
$d.Paragraphs.Item($anchorIndex + 2).Range.Delete()
$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()
